# testing_taks2.xlsx — add Task2 "C" (average) block + fix Task B "multiply" row.
#
# Summary of changes (see commit message: "add task2 def test_average and
# fix(main): def multi"):
#   1. Legend table (I:K) gains a "C" row (row 6) describing the new block.
#   2. Row 20 (task B, "7.B" / "другой массив с числами") G column result
#      is corrected from TypeError to Failed.
#   3. A brand-new "Testing task2 Array def avarage" block is appended
#      (rows 23-33): a merged title row, a header row, and 9 data rows
#      (1.C .. 9.C).
#   4. Column G is widened to fit the new, longer "ZeroDivisionError" value.
#   5. Selection cursor moves to J14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Helper pattern used throughout: copy a well-formatted source range's
# *formatting* onto a destination range (values are written separately),
# so every new cell matches the look of its sibling rows/tables exactly.
# ---------------------------------------------------------------------

# --- 1. Legend table: add the "C" row under A/B rows ------------------
$ws.Range("B7:C7").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("K5").Copy()
$ws.Range("K6").PasteSpecial(-4122)

$ws.Range("I6").Value = "C"
$ws.Range("J6").Value = "Среднее арифметическое"
$ws.Range("K6").Value = "на разных входных данных"

# --- 2. Fix Task B row 20 ("7.B") result: TypeError -> Failed ---------
$ws.Range("G20").Value = "Failed"

# --- 3. New "Testing task2 Array def avarage" block --------------------

# 3a. Title row (merged A23:B23), formatted like the other section titles.
$ws.Range("A12:B12").Copy()
$ws.Range("A23:B23").PasteSpecial(-4122)
$ws.Range("A23:B23").Merge()
$ws.Range("A23").Value = "Testing task2 Array def avarage"

# 3b. Header row (A24:G24), formatted like the other header rows.
$ws.Range("A13:G13").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$ws.Range("A24").Value = "ID"
$ws.Range("B24").Value = "name "
$ws.Range("C24").Value = "priority"
$ws.Range("D24").Value = "pred condition"
$ws.Range("E24").Value = "input data"
$ws.Range("F24").Value = "expected"
$ws.Range("G24").Value = "result"

# 3c. Data rows 25-31 follow the same formatting as Task A's 3-9
#     (only IDs/expected values change: sum -> average).
$ws.Range("A3:G9").Copy()
$ws.Range("A25:G31").PasteSpecial(-4122)

# Row 32 is an extra case not present in Task A; base its format on row 8
# (same "wrong type" styling: priority 2, highlighted name/priority cells).
$ws.Range("A8:G8").Copy()
$ws.Range("A32").PasteSpecial(-4122)

# Row 33 mirrors Task A's row 10 ("False вместо элемента").
$ws.Range("A10:G10").Copy()
$ws.Range("A33").PasteSpecial(-4122)

# --- 3d. Values for rows 25-33 -----------------------------------------

# 1.C - целые числа
$ws.Range("A25").Value = "1.C"
$ws.Range("B25").Value = "целые числа"
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = "array"
$ws.Range("E25").Value = "[1,2,3,4,5]"
$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 3

# 2.C - вещественные числа
$ws.Range("A26").Value = "2.C"
$ws.Range("B26").Value = "вещественные числа"
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = "array"
$ws.Range("E26").Value = "[1.2,2,3.8,4,5]"
$ws.Range("F26").Value = 3.2
$ws.Range("G26").Value = 3.2

# 3.C - отрицательные числа
$ws.Range("A27").Value = "3.C"
$ws.Range("B27").Value = "отрицательные числа"
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = "array"
$ws.Range("E27").Value = "[-1,-2,-3,-4,-5]"
$ws.Range("F27").Value = -3
$ws.Range("G27").Value = -3

# 4.C - разные числа
$ws.Range("A28").Value = "4.C"
$ws.Range("B28").Value = "разные числа"
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = "array"
$ws.Range("E28").Value = "[1,-2,-3,4,5]"
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1

# 5.C - пустой массив -> ZeroDivisionError
$ws.Range("A29").Value = "5.C"
$ws.Range("B29").Value = "пустой массив"
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = "array"
$ws.Range("E29").Value = "[]"
$ws.Range("F29").Value = "ZeroDivisionError"
$ws.Range("G29").Value = "ZeroDivisionError"

# 6.C - разные типы
$ws.Range("A30").Value = "6.C"
$ws.Range("B30").Value = "разные типы"
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = "array"
$ws.Range("E30").Value = '[1,"2",3,"4",5]'
$ws.Range("F30").Value = "TypeError"
$ws.Range("G30").Value = "TypeError"

# 7.C - другой массив с числами
$ws.Range("A31").Value = "7.C"
$ws.Range("B31").Value = "другой массив с числами"
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = "array"
$ws.Range("E31").Value = "[1,2,[1,2],4,5]"
$ws.Range("F31").Value = "TypeError"
$ws.Range("G31").Value = "TypeError"

# 8.C - разные типы (extra case)
$ws.Range("A32").Value = "8.C"
$ws.Range("B32").Value = "разные типы"
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = "array"
$ws.Range("E32").Value = '["", 1, {1}, "", None]'
$ws.Range("F32").Value = "TypeError"
$ws.Range("G32").Value = "TypeError"

# 9.C - False вместо элемента
$ws.Range("A33").Value = "9.C"
$ws.Range("B33").Value = "False вместо элемента"
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = "array"
$ws.Range("E33").Value = "[1,2, False,4,5]"
$ws.Range("F33").Value = "TypeError"
$ws.Range("G33").Value = "Failed"

# --- 4. Column G is now wider (fits "ZeroDivisionError") ---------------
$ws.Columns.Item(7).ColumnWidth = 15.59

# --- 5. Move the selection cursor like the author's last edit ----------
$ws.Activate()
$ws.Range("J14").Select()
